$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.297.53"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "2.441.68"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.12"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.11"
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -2.65%  "
$ws.Range("E9").Value = "  -2.00%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.18"
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.345"
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.46"
$ws.Range("E13").Value = "  -2.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000172"
$ws.Range("E14").Value = "  -3.94%  "
$ws.Range("D15").Value = "2.885.62"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").Value = "62.134.92"
$ws.Range("E16").Value = "  -1.78%  "
$ws.Range("D17").Value = "2.448.82"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.62"
$ws.Range("E18").Value = "  -4.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.66"
$ws.Range("E19").Value = "  -4.23%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "319.30"
$ws.Range("E20").Value = "  -3.59%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.10"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("E22").Value = "  -3.13%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.71"
$ws.Range("E24").Value = "  +3.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.76"
$ws.Range("E25").Value = "  -2.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "633.82"
$ws.Range("E26").Value = "  -5.24%  "
$ws.Range("D27").Value = "2.560.15"
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0944"
$ws.Range("E29").Value = "  -5.88%  "
$ws.Range("E30").Value = "  -5.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.81"
$ws.Range("E31").Value = "  -4.14%  "
$ws.Range("E32").Value = "  -4.24%  "
$ws.Range("E33").Value = "  -3.85%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  -4.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.60"
$ws.Range("E36").Value = "  -4.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "150.40"
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.362"
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("E40").Value = "  -5.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.67"
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("E42").Value = "  -4.02%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +1.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "151.02"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.91"
$ws.Range("E49").Value = "  -5.13%  "
$ws.Range("E50").Value = "  -3.18%  "
$ws.Range("E51").Value = "  -2.46%  "
